# Fruta / hortaliza, semanal
# Weekly price-series update: a new week's record is inserted at the top of
# the "Zapallo italiano" price history (row 481), pushing every existing
# record down by one row (old 481-507 become new 482-508).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 481, shifting the rest down.
$ws.Rows(481).Insert()

# Populate the newly inserted row with this week's data.
$ws.Cells.Item(481, 1).Value  = 4
$ws.Cells.Item(481, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(481, 3).Value  = "Los Lagos"
$ws.Cells.Item(481, 4).Value  = 45267
$ws.Cells.Item(481, 5).Value  = 10
$ws.Cells.Item(481, 6).Value  = 100112032
$ws.Cells.Item(481, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(481, 8).Value  = "Sin especificar"
$ws.Cells.Item(481, 9).Value  = "Primera"
$ws.Cells.Item(481, 10).Value = 150
$ws.Cells.Item(481, 11).Value = 16000
$ws.Cells.Item(481, 12).Value = 16000
$ws.Cells.Item(481, 13).Value = 16000
$ws.Cells.Item(481, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(481, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(481, 16).Value = 320
$ws.Cells.Item(481, 17).Value = 50
$ws.Cells.Item(481, 18).Value = "Hortaliza"
